$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: 893×2=1786 -> 967×4=3868
$cell = $t.Cell(1,1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "893×2=1786") {
    Write-Host "WARNING: Row 1 Col 1 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "967×4=3868"

# Row 1, Col 2: 826×9=7434 -> 794×9=7146
$cell = $t.Cell(1,2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "826×9=7434") {
    Write-Host "WARNING: Row 1 Col 2 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "794×9=7146"

# Row 1, Col 3: 799×7=5593 -> 674×2=1348
$cell = $t.Cell(1,3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "799×7=5593") {
    Write-Host "WARNING: Row 1 Col 3 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "674×2=1348"

# Row 1, Col 4: 995×6=5970 -> 206×7=1442
$cell = $t.Cell(1,4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "995×6=5970") {
    Write-Host "WARNING: Row 1 Col 4 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "206×7=1442"

# Row 1, Col 5: 967×8=7736 -> 543×8=4344
$cell = $t.Cell(1,5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "967×8=7736") {
    Write-Host "WARNING: Row 1 Col 5 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "543×8=4344"

# Row 5, Col 1: 285×8=2280 -> 452×6=2712
$cell = $t.Cell(5,1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "285×8=2280") {
    Write-Host "WARNING: Row 5 Col 1 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "452×6=2712"

# Row 5, Col 2: 383×2=766 -> 508×8=4064
$cell = $t.Cell(5,2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "383×2=766") {
    Write-Host "WARNING: Row 5 Col 2 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "508×8=4064"

# Row 5, Col 3: 824×7=5768 -> 171×7=1197
$cell = $t.Cell(5,3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "824×7=5768") {
    Write-Host "WARNING: Row 5 Col 3 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "171×7=1197"

# Row 5, Col 4: 292×2=584 -> 453×3=1359
$cell = $t.Cell(5,4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "292×2=584") {
    Write-Host "WARNING: Row 5 Col 4 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "453×3=1359"

# Row 5, Col 5: 549×6=3294 -> 725×9=6525
$cell = $t.Cell(5,5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "549×6=3294") {
    Write-Host "WARNING: Row 5 Col 5 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "725×9=6525"

# Row 10, Col 1: 741×3=2223 -> 977×3=2931
$cell = $t.Cell(10,1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "741×3=2223") {
    Write-Host "WARNING: Row 10 Col 1 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "977×3=2931"

# Row 10, Col 2: 722×2=1444 -> 367×4=1468
$cell = $t.Cell(10,2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "722×2=1444") {
    Write-Host "WARNING: Row 10 Col 2 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "367×4=1468"

# Row 10, Col 3: 466×9=4194 -> 598×7=4186
$cell = $t.Cell(10,3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "466×9=4194") {
    Write-Host "WARNING: Row 10 Col 3 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "598×7=4186"

# Row 10, Col 4: 919×8=7352 -> 623×5=3115
$cell = $t.Cell(10,4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "919×8=7352") {
    Write-Host "WARNING: Row 10 Col 4 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "623×5=3115"

# Row 10, Col 5: 549×6=3294 -> 610×5=3050
$cell = $t.Cell(10,5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "549×6=3294") {
    Write-Host "WARNING: Row 10 Col 5 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "610×5=3050"

# Row 15, Col 1: 548×9=4932 -> 262×2=524
$cell = $t.Cell(15,1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "548×9=4932") {
    Write-Host "WARNING: Row 15 Col 1 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "262×2=524"

# Row 15, Col 2: 109×5=545 -> 619×8=4952
$cell = $t.Cell(15,2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "109×5=545") {
    Write-Host "WARNING: Row 15 Col 2 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "619×8=4952"

# Row 15, Col 3: 149×5=745 -> 122×8=976
$cell = $t.Cell(15,3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "149×5=745") {
    Write-Host "WARNING: Row 15 Col 3 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "122×8=976"

# Row 15, Col 4: 396×3=1188 -> 982×3=2946
$cell = $t.Cell(15,4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "396×3=1188") {
    Write-Host "WARNING: Row 15 Col 4 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "982×3=2946"

# Row 15, Col 5: 820×3=2460 -> 350×8=2800
$cell = $t.Cell(15,5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "820×3=2460") {
    Write-Host "WARNING: Row 15 Col 5 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "350×8=2800"

# Row 20, Col 1: 333×6=1998 -> 678×3=2034
$cell = $t.Cell(20,1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "333×6=1998") {
    Write-Host "WARNING: Row 20 Col 1 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "678×3=2034"

# Row 20, Col 2: 980×2=1960 -> 756×2=1512
$cell = $t.Cell(20,2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "980×2=1960") {
    Write-Host "WARNING: Row 20 Col 2 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "756×2=1512"

# Row 20, Col 3: 823×9=7407 -> 333×4=1332
$cell = $t.Cell(20,3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "823×9=7407") {
    Write-Host "WARNING: Row 20 Col 3 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "333×4=1332"

# Row 20, Col 4: 101×7=707 -> 562×5=2810
$cell = $t.Cell(20,4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "101×7=707") {
    Write-Host "WARNING: Row 20 Col 4 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "562×5=2810"

# Row 20, Col 5: 557×2=1114 -> 349×3=1047
$cell = $t.Cell(20,5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "557×2=1114") {
    Write-Host "WARNING: Row 20 Col 5 unexpected content: [$($cell.Range.Text)]"
}
$cell.Range.Text = "349×3=1047"

Write-Host "Done."